$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A45:A47 with the shared string already used in A44 ("Alaska Washington
# and Oregon fishermen licensed in California") and give the rows the same
# row height / styling as A44 so the merged label reads across rows 44-47.
$labelText = $ws.Range("A44").Value2

$ws.Range("A45").Value = $labelText
$ws.Range("A46").Value = $labelText
$ws.Range("A47").Value = $labelText

$ws.Range("A45").Style = $ws.Range("A44").Style
$ws.Range("A46").Style = $ws.Range("A44").Style
$ws.Range("A47").Style = $ws.Range("A44").Style

$ws.Rows.Item(45).RowHeight = 51
$ws.Rows.Item(46).RowHeight = 51
$ws.Rows.Item(47).RowHeight = 51

# Update the view state to match the saved selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("A47").Select()
